$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('C2').Value = 'Загрузка…'

$ws.Range('C3').Value = 'О нет! В лаборатории произошла авария и ты застряла в огне!'

$ws.Range('C4').Value = 'Твоя армия обученных обезьянок сбежала…'
$ws.Range('C4').Font.Name = "Calibri"

$ws.Range('C5').Value = '…все, кроме одной.'
$ws.Range('C5').Font.Name = "Calibri"

$ws.Range('C6').Value = 'одной обезьянки.'
$ws.Range('C6').Font.Name = "Calibri"

$ws.Range('C7').Value = 'дай ей инструкции,'
$ws.Range('C7').Font.Name = "Calibri"

$ws.Range('C8').Value = 'чтобы спастись!'
$ws.Range('C8').Font.Name = "Calibri"

$ws.Range('C9').Value = 'Авторы:'
$ws.Range('C9').Font.Name = "Calibri"

$ws.Range('C10').Value = 'Код:'
$ws.Range('C10').Font.Name = "Calibri"

$ws.Range('C11').Value = 'Anders Antila'

$ws.Range('C12').Value = 'Иллюстрации:'
$ws.Range('C12').Font.Name = "Calibri"

$ws.Range('C13').Value = 'erik Nahkala & laura kantti'

$ws.Range('C14').Value = 'Музыка:'
$ws.Range('C14').Font.Name = "Calibri"

$ws.Range('C15').Value = 'RaxL snaxel'

$ws.Range('C16').Value = 'Побег Обезьянок'
$ws.Range('C16').Font.Name = "Calibri"

$ws.Range('C17').Value = 'Кликните, чтобы начать!'
$ws.Range('C17').Font.Name = "Calibri"

$ws.Range('C18').Value = '*или удивительная история о том, как я спаслась от пожара и катастрофы благодаря надежному и верному другу-обезьянке)'
$ws.Range('C18').Font.Name = "Calibri"

$ws.Range('C19').Value = 'У тебя получилось!'
$ws.Range('C19').Font.Name = "Calibri"

$ws.Range('C20').Value = 'Назад'
$ws.Range('C20').Font.Name = "Calibri"

$ws.Range('C21').Value = 'авторы'
$ws.Range('C21').Font.Name = "Calibri"

$ws.Range('C22').Value = ' Уровень 1'
$ws.Range('C22').Font.Name = "Calibri"

$ws.Range('C23').Value = ' Уровень 2'
$ws.Range('C23').Font.Name = "Calibri"

$ws.Range('C24').Value = ' Уровень 3'
$ws.Range('C24').Font.Name = "Calibri"

$ws.Range('C25').Value = ' Уровень 4'
$ws.Range('C25').Font.Name = "Calibri"

$ws.Range('C26').Value = ' Уровень 5'
$ws.Range('C26').Font.Name = "Calibri"

$ws.Range('C27').Value = ' Уровень 6'
$ws.Range('C27').Font.Name = "Calibri"

$ws.Range('C28').Value = ' Уровень 7'
$ws.Range('C28').Font.Name = "Calibri"

$ws.Range('C29').Value = ' Уровень 8'
$ws.Range('C29').Font.Name = "Calibri"

$ws.Range('C30').Value = 'меню'
$ws.Range('C30').Font.Name = "Calibri"

$ws.Range('C31').Value = 'выкл. звук'
$ws.Range('C31').Font.Name = "Calibri"

$ws.Range('C32').Value = 'дальше'
$ws.Range('C32').Font.Name = "Calibri"

$ws.Range('C33').Value = 'играть'
$ws.Range('C33').Font.Name = "Calibri"

$ws.Range('C34').Value = 'вкл. звук'
$ws.Range('C34').Font.Name = "Calibri"

$ws.Range('C35').Value = 'язык'
$ws.Range('C35').Font.Name = "Calibri"

$ws.Range('D18').Value = '( * o la fantastica storia di quella volta in cui mi sono stato salvato dalle fiamme e disavventure grazie all''aiuto di una fedele e leale scimmietta)'

$ws.Rows.Item(2).RowHeight = 16
$ws.Rows.Item(3).RowHeight = 32
$ws.Rows.Item(4).RowHeight = 16
$ws.Rows.Item(5).RowHeight = 16
$ws.Rows.Item(6).RowHeight = 16
$ws.Rows.Item(7).RowHeight = 16
$ws.Rows.Item(8).RowHeight = 16
$ws.Rows.Item(9).RowHeight = 16
$ws.Rows.Item(10).RowHeight = 16
$ws.Rows.Item(11).RowHeight = 16
$ws.Rows.Item(12).RowHeight = 16
$ws.Rows.Item(13).RowHeight = 16
$ws.Rows.Item(14).RowHeight = 16
$ws.Rows.Item(15).RowHeight = 16
$ws.Rows.Item(16).RowHeight = 16
$ws.Rows.Item(17).RowHeight = 16
$ws.Rows.Item(18).RowHeight = 48
$ws.Rows.Item(19).RowHeight = 16
$ws.Rows.Item(20).RowHeight = 16
$ws.Rows.Item(21).RowHeight = 16
$ws.Rows.Item(22).RowHeight = 16
$ws.Rows.Item(23).RowHeight = 16
$ws.Rows.Item(24).RowHeight = 16
$ws.Rows.Item(25).RowHeight = 16
$ws.Rows.Item(26).RowHeight = 16
$ws.Rows.Item(27).RowHeight = 16
$ws.Rows.Item(28).RowHeight = 16
$ws.Rows.Item(29).RowHeight = 16
$ws.Rows.Item(30).RowHeight = 16
$ws.Rows.Item(31).RowHeight = 16
$ws.Rows.Item(32).RowHeight = 16
$ws.Rows.Item(33).RowHeight = 16
$ws.Rows.Item(34).RowHeight = 16
$ws.Rows.Item(35).RowHeight = 16

$ws.Range("C38").Select()
